$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell AT1 for semana epidemiologica 43 (text, matching style of existing week-number headers)
$ws.Range("AT1").NumberFormat = "@"
$ws.Range("AT1").Value = "43"

# Numeric data updates / additions for week 43 column (AT) and a handful of corrected values
$ws.Range("AT2").Value = 0
$ws.Range("AT3").Value = 0
$ws.Range("AT5").Value = 0
$ws.Range("AS6").Value = 26
$ws.Range("AT6").Value = 20
$ws.Range("AT7").Value = 1
$ws.Range("AT8").Value = 8
$ws.Range("AT9").Value = 0
$ws.Range("AT10").Value = 0
$ws.Range("AT12").Value = 0
$ws.Range("AT13").Value = 0
$ws.Range("AT14").Value = 0
$ws.Range("AT15").Value = 0
$ws.Range("AT16").Value = 0
$ws.Range("AT17").Value = 0
$ws.Range("AT18").Value = 0
$ws.Range("AT19").Value = 0
$ws.Range("AT23").Value = 0
$ws.Range("AT25").Value = 8
$ws.Range("AT26").Value = 0
$ws.Range("AT28").Value = 5
$ws.Range("AT29").Value = 6
$ws.Range("L29").Value = 3
$ws.Range("AO30").Value = 1
$ws.Range("AT30").Value = 2
$ws.Range("AT31").Value = 0
$ws.Range("AR34").Value = 3
$ws.Range("AS34").Value = 0
$ws.Range("AT34").Value = 0
$ws.Range("AS35").Value = 4
$ws.Range("AT35").Value = 6
$ws.Range("AH36").Value = 0
$ws.Range("AJ36").Value = 0
$ws.Range("AK36").Value = 2
$ws.Range("AL36").Value = 1
$ws.Range("AP36").Value = 0
$ws.Range("AT36").Value = 0
$ws.Range("W36").Value = 1
$ws.Range("Y36").Value = 0
$ws.Range("Z36").Value = 0
$ws.Range("AT37").Value = 0
$ws.Range("AT38").Value = 0
$ws.Range("AT42").Value = 0
$ws.Range("AT43").Value = 0
$ws.Range("AT44").Value = 0
$ws.Range("AT45").Value = 0
$ws.Range("AT46").Value = 0
$ws.Range("AT47").Value = 0
$ws.Range("AT48").Value = 0
$ws.Range("AT49").Value = 0
$ws.Range("AT50").Value = 0
$ws.Range("AT51").Value = 0
$ws.Range("AT53").Value = 0
$ws.Range("AT54").Value = 0
$ws.Range("AT55").Value = 0
$ws.Range("AT56").Value = 0
$ws.Range("AT57").Value = 0
$ws.Range("AT58").Value = 0
